$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "quantity" in C1
$ws.Range("C1").Value = "quantity"

# Fill column C with quantity values (C2:C7 = 2)
$ws.Range("C2:C7").Value = 2

# Fix B7 value from 0 to 2
$ws.Range("B7").Value = 2

# Auto-fit column B to content width (header "burning time")
$ws.Columns("B:B").AutoFit()

# Selection per diff (activeCell C9)
$ws.Range("C9").Select()
